$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks (email contact links for the old dataset).
# We'll re-create the ones that survive (rows 2-4) further down so that the
# relationship ids / styles come out clean.
$ws.Hyperlinks.Delete()

# --- Row 2: puntosamarillo / mario hernandez -----------------------------
$ws.Range("A2").Value = "puntosamarillo"
$ws.Range("B2").Value = "puntosamarillo.com"
$ws.Range("C2").Value = "pinturas"
$ws.Range("D2").Value = "mario"
$ws.Range("E2").Value = "hernandez"
$ws.Range("F2").Value = 3203525634
$ws.Range("G2").Value = "mario@puntosamarillo.com"
$ws.Range("H2").Value = "precualification"

# --- Row 3: osomazorca / oscar agudelo -----------------------------------
$ws.Range("A3").Value = "osomazorca"
$ws.Range("B3").Value = "osomazorca.com"
$ws.Range("C3").Value = "comidas"
$ws.Range("D3").Value = "oscar"
$ws.Range("E3").Value = "agudelo"
$ws.Range("F3").Value = 3112324563
$ws.Range("G3").Value = "oscar@osomazorca.com"
$ws.Range("H3").Value = "precualification"

# --- Row 4: toshiba / mariano carreno -------------------------------------
$ws.Range("A4").Value = "toshiba"
$ws.Range("B4").Value = "toshiba.com"
$ws.Range("C4").Value = "computers"
$ws.Range("D4").Value = "mariano"
$ws.Range("E4").Value = "carreno"
$ws.Range("F4").Value = 3102343267
$ws.Range("G4").Value = "mariano@toshiba.com"
$ws.Range("H4").Value = "precualification"

# Re-add the mailto hyperlinks for the contact e-mails, then restore the
# plain (non-hyperlink) font so the cell keeps its original style.
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:mario@puntosamarillo.com", "", "", "mario@puntosamarillo.com")
$ws.Range("G2").Font.Name = "Calibri"
$ws.Range("G2").Font.Size = 11
$ws.Range("G2").Font.Color = 0
$ws.Range("G2").Font.Underline = $false

$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:oscar@osomazorca.com", "", "", "oscar@osomazorca.com")
$ws.Range("G3").Font.Name = "Calibri"
$ws.Range("G3").Font.Size = 11
$ws.Range("G3").Font.Color = 0
$ws.Range("G3").Font.Underline = $false

$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:mariano@toshiba.com", "", "", "mariano@toshiba.com")
$ws.Range("G4").Font.Name = "Calibri"
$ws.Range("G4").Font.Size = 11
$ws.Range("G4").Font.Color = 0
$ws.Range("G4").Font.Underline = $false

# --- Rows 5-15: wipe out all the leftover sample/test data ---------------
$ws.Range("A5:H15").ClearContents()

# Move the active selection to C8 (matches the saved selection state).
$ws.Range("C8").Select()
